# End of Session: H-044 Context Helper System & Maintenance Expansion
# Update floor plan column widths on the "Column Definitions" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Column Definitions")

# Adjust the "width" values for code / bedroom_count / bathroom_count rows.
$ws.Range("C2").Value = 140
$ws.Range("C3").Value = 60
$ws.Range("C4").Value = 150

# Move the active selection as left by the author at end of session.
$ws.Activate()
$ws.Range("F14").Select()
